# Apply the "Deploying to gh-pages" content update to
# StructureDefinition-communication-measure.xlsx:
#   - Metadata sheet: Version 5.0.0 -> 6.0.0, Date bump, Publisher filled in
#     as "Alvearie Team", the second (duplicate) "Contact" row is replaced
#     by a new "Jurisdiction" / "United States of America" row, and the
#     leftover duplicate "Contact" row is removed entirely.
#   - Elements sheet: the root Extension row's Short/Definition text is
#     updated to match the StructureDefinition's own Title/Description.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes the new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The second, now-redundant "Contact" row (row 11) is deleted outright,
# shifting every following row up by one.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition now mirror the StructureDefinition's
# own Title / Description.
$elements.Range("K2").Value = "Communication Measure"
$elements.Range("L2").Value = "Reference to the measure that resulted in the communication"
